$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: change A2 from "Saurab@gmail.com" to "Amol@gmail.com" and hyperlink it ---
$ws.Range("A2").Value = "Amol@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:Amol@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

# --- New "Status" column (D) ---
$ws.Range("D1").Value = "Status"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D2").Value = "Pass"

# --- Row 3: new username/password pair ---
$ws.Range("A3").Value = "Chandu@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:Chandu@gmail.com")
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("B3").Value = "pass1123"

# --- Row 4: new username/password pair ---
$ws.Range("A4").Value = "Azhar@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:Azhar@gmail.com")
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("B4").Value = "test345"

# --- Touch E4 so the sheet's used range/dimension extends to E4 (matches target) ---
$ws.Range("E4").Font.Bold = $false

# --- Final selection lands on E10, matching the saved view state ---
[void]$ws.Range("E10").Select()
